# This script applies a weekly data update to the "Naranja" price sheet:
# 7 new price rows (dated 2021-11-24, serial 44524) are inserted right
# before the current row 277, pushing every subsequent row down by 7
# (old row 277 -> new row 284, ..., old row 365 -> new row 372).
# The rest of the sheet's data is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 blank rows at 277-283, shifting rows 277:365 down to 284:372.
$ws.Range("A277:T283").Insert()

# Common / repeated values for every detail row on this sheet.
$mercadoId   = 2
$mercado     = "Comercializadora del Agro de Limarí"
$region      = "Coquimbo"
$codreg      = 4
$tipo        = "Fruta"
$productoId  = 100102
$producto    = "Cítricos"
$categoriaId = 100102005
$categoria   = "Naranja"
$unidad      = "$/bins (400 kilos)"
$origen      = "Provincia de Limarí"
$kgUnidad    = 400
$fecha       = 44524

function Set-PriceRow {
    param($row, $variedad, $calidad, $volumen, $min, $max, $prom, $precioKg)

    $ws.Range("A$row").Value = $mercadoId
    $ws.Range("B$row").Value = $mercado
    $ws.Range("C$row").Value = $region
    $ws.Range("D$row").Value = $fecha
    $ws.Range("E$row").Value = $codreg
    $ws.Range("F$row").Value = $tipo
    $ws.Range("G$row").Value = $productoId
    $ws.Range("H$row").Value = $producto
    $ws.Range("I$row").Value = $categoriaId
    $ws.Range("J$row").Value = $categoria
    $ws.Range("K$row").Value = $variedad
    $ws.Range("L$row").Value = $calidad
    $ws.Range("M$row").Value = $volumen
    $ws.Range("N$row").Value = $min
    $ws.Range("O$row").Value = $max
    $ws.Range("P$row").Value = $prom
    $ws.Range("Q$row").Value = $unidad
    $ws.Range("R$row").Value = $origen
    $ws.Range("S$row").Value = $precioKg
    $ws.Range("T$row").Value = $kgUnidad
}

Set-PriceRow 277 "Cara cara"  "Primera" 20 175000 180000 177500 444
Set-PriceRow 278 "Cara cara"  "Segunda" 20 145000 150000 147500 369
Set-PriceRow 279 "Lane Late"  "Primera" 20 165000 170000 167500 419
Set-PriceRow 280 "Lane Late"  "Segunda" 16 145000 150000 147500 369
Set-PriceRow 281 "Navel Late" "Primera" 16 165000 170000 167500 419
Set-PriceRow 282 "Navel Late" "Segunda" 16 145000 150000 147500 369
Set-PriceRow 283 "Valencia"   "Primera" 20 175000 180000 177500 444
